$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Intercept)
$ws.Range("B2").Value = "0.237***`n (0.087)"
$ws.Range("D2").Value = "0.355***`n (0.060)"
$ws.Range("F2").Value = "0.521***`n (0.047)"

# Row 3 (Mining)
$ws.Range("B3").Value = "0.010`n (0.019)"
$ws.Range("C3").Value = "0.049`n (0.055)"
$ws.Range("D3").Value = "0.008`n (0.012)"
$ws.Range("E3").Value = "0.005`n (0.034)"
$ws.Range("F3").Value = "-0.017`n (0.010)"
$ws.Range("G3").Value = "0.003`n (0.023)"

# Row 4 (Partisanship)
$ws.Range("B4").Value = "-0.385***`n (0.099)"
$ws.Range("C4").Value = "-0.545**`n (0.214)"
$ws.Range("D4").Value = "0.504***`n (0.067)"
$ws.Range("E4").Value = "0.294**`n (0.134)"
$ws.Range("F4").Value = "0.376***`n (0.052)"
$ws.Range("G4").Value = "0.109`n (0.091)"

# Row 5 (Deregulated)
$ws.Range("B5").Value = "0.291*`n (0.169)"
$ws.Range("D5").Value = "0.078`n (0.113)"
$ws.Range("F5").Value = "0.179*`n (0.093)"

# Row 6 (R-squared)
$ws.Range("B6").Value = 0.2962603545018835
$ws.Range("C6").Value = 0.130089143551553
$ws.Range("D6").Value = 0.5616280652897079
$ws.Range("E6").Value = 0.1011422130944299
$ws.Range("F6").Value = 0.4634153834670367
$ws.Range("G6").Value = 0.02731738325835731

# Row 7 (N)
$ws.Range("B7").Value = 63
$ws.Range("C7").Value = 63
$ws.Range("D7").Value = 55
$ws.Range("E7").Value = 55
$ws.Range("F7").Value = 64
$ws.Range("G7").Value = 64
